$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting existing rows 61-70 down to 62-71
$ws.Rows.Item(61).Insert()

# Populate the new row 61 with the new weekly data point
$ws.Range("A61").Value = 4
$ws.Range("B61").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C61").Value = "Los Lagos"
$ws.Range("D61").Value = 44663
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = 100112031
$ws.Range("G61").Value = "Poroto verde"
$ws.Range("H61").Value = "Magnum"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 60
$ws.Range("K61").Value = 28000
$ws.Range("L61").Value = 28000
$ws.Range("M61").Value = 28000
$ws.Range("N61").Value = "$/saco 25 kilos"
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 1120
$ws.Range("Q61").Value = 25
$ws.Range("R61").Value = "Hortaliza"
